$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

$xlCenter = -4108

# --- New "target" / "suoritettu(%)" rows (written first so the shared-
# string table picks up "target"/"suoritettu(%)" before the longer log
# entries below, matching the saved file's string order). ---
$ws.Cells.Item(87, 1).Value = "target"
$ws.Cells.Item(87, 1).HorizontalAlignment = $xlCenter

$ws.Cells.Item(88, 1).Value = "suoritettu(%)"
$ws.Cells.Item(88, 1).HorizontalAlignment = $xlCenter

# New log entries (rows 80 and 81).
$ws.Cells.Item(80, 2).Value = 1
$ws.Cells.Item(80, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(80, 2).VerticalAlignment = $xlCenter
$ws.Cells.Item(80, 3).Value = "errorMessage komponentti muutettu Notification yleiskomponentiksi  (error, success, action jaot), reducer muokattu"
$ws.Cells.Item(80, 4).Value = "client"

$ws.Cells.Item(81, 2).Value = 1
$ws.Cells.Item(81, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(81, 2).VerticalAlignment = $xlCenter
$ws.Cells.Item(81, 3).Value = "Notification otettu käyttöön Login ja Sign Up komponenteissa laajemmin, muissa error osio käytössä normaalisti"
$ws.Cells.Item(81, 4).Value = "client"

# Totals row: extend the sum to include the new rows.
$ws.Range("B86").Formula = "=SUM(B2:B81)"

# Finish the new target / completion rows.
$ws.Cells.Item(87, 2).Value = 175
$ws.Cells.Item(87, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(87, 2).VerticalAlignment = $xlCenter

$ws.Range("B88").Formula = "=B86/B87*100"
$ws.Cells.Item(88, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(88, 2).VerticalAlignment = $xlCenter

# Widen column A slightly to fit the new "suoritettu(%)" label
# (target raw width 14.7109375; the COM width model here snaps to
# sixth-character increments, so 13.8 is the closest reachable input).
$ws.Columns.Item(1).ColumnWidth = 13.8

# Update the active selection to mirror the edited cell.
$ws.Range("B78").Select()
